# Update cryptos list cell values to match the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.717.30"
$ws.Range("E2").Value = "  -1.24%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.639.12"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "579.31"
$ws.Range("E5").Value = "  +0.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "155.23"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.621"
$ws.Range("E8").Value = "  -4.16%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.634.90"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  -4.18%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.79"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  +0.92%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.38"
$ws.Range("E14").Value = "  -0.88%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.114.71"
$ws.Range("E15").Value = "  +0.28%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000183"
$ws.Range("E16").Value = "  -2.30%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "63.676.67"
$ws.Range("E17").Value = "  -1.01%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.638.90"
$ws.Range("E18").Value = "  -0.41%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.12"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.65"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("E21").Value = "  -3.22%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "344.38"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -4.55%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "602.90"
$ws.Range("E27").Value = "  +6.47%  "
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("E29").Value = "  +1.36%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.12"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.161"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("E33").Value = "  -0.75%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.73"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("E35").Value = "  -1.60%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.45"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("E38").Value = "  +0.04%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.66"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("E40").Value = "  -2.52%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "150.80"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  +3.28%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "41.91"
$ws.Range("E44").Value = "  -0.69%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "161.08"
$ws.Range("E45").Value = "  +1.63%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "24.14"
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("E49").Value = "  -0.70%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0997"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("E51").Value = "  -1.40%  "
